# Apply the spring-2020 vote updates to the "lablist" sheet.
# These are the only real user-entered changes; every Y/Z/AA/V/T81-86/S84-86
# cell in the sheet is a formula and recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("lablist")

# Row 31 ("Quick and Easy" info session?) — U vote of 2 removed
$ws.Range("U31").ClearContents()

# Row 53 — add a 0.5 T-column vote
$ws.Range("T53").Value = 0.5

# Row 54 — add a 1 T-column vote
$ws.Range("T54").Value = 1

# Row 64 — add a 1 T-column vote
$ws.Range("T64").Value = 1

# Rows 67-70 — remove the U-column vote of 1
$ws.Range("U67").ClearContents()
$ws.Range("U68").ClearContents()
$ws.Range("U69").ClearContents()
$ws.Range("U70").ClearContents()

$ws.Calculate()
